$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-row "week" data which gets shifted down by one row
# (D = Fecha, J = Volumen, K = Precio minimo, L = Precio maximo, M = Precio promedio,
#  O = Origen, P = Precio $/Kg)
$cols = @("D","J","K","L","M","O","P")

# 1) Capture the current ("before") values for rows 68..135 for the shifting columns,
#    plus the full row 135 (all columns) which will become the new row 136.
$old = @{}
for ($r = 68; $r -le 135; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $old[$r] = $rowVals
}

$fullCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
$oldRow135 = @{}
foreach ($c in $fullCols) {
    $oldRow135[$c] = $ws.Range("${c}135").Value2
}

# 2) Create the new row 136 as an exact copy of the old row 135.
foreach ($c in $fullCols) {
    $ws.Range("${c}136").Value2 = $oldRow135[$c]
}
$ws.Range("D136").NumberFormat = $ws.Range("D135").NumberFormat

# 3) Shift the data down: new row r (for r = 135 down to 69) takes the old values of row r-1.
for ($r = 135; $r -ge 69; $r--) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $old[$r - 1][$c]
    }
}

# 4) Row 68 gets the brand-new week's values.
$ws.Range("D68").Value2 = 44447
$ws.Range("J68").Value2 = 45
$ws.Range("K68").Value2 = 10000
$ws.Range("L68").Value2 = 10000
$ws.Range("M68").Value2 = 10000
$ws.Range("O68").Value2 = "Región de Arica y Parinacota"
$ws.Range("P68").Value2 = 167
